# Task SummaryWk11.xlsx -- "moved local copies to repo"
# Fill in the week's task-summary data (author + week number, the four
# task rows, the cumulative-total line) and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: author name + week number
$ws.Range("C1").Value = "Jesse Hare"
$ws.Range("E1").Value = 11

# Stage column for the four task rows (3-6): all "Project Build"
$ws.Range("A3").Value = "Project Build"
$ws.Range("A4").Value = "Project Build"
$ws.Range("A5").Value = "Project Build"
$ws.Range("A6").Value = "Project Build"

# Totals row: cumulative total label (literal total replaces old SUM formula)
$ws.Range("A14").Value = "Cumulative Total: 220"

# Task descriptions for rows 3-6
$ws.Range("B3").Value = "Implementation of final feature test"
$ws.Range("B4").Value = "Unit Testing"
$ws.Range("B5").Value = "Integration Testing (with front end)"
$ws.Range("B6").Value = "Finish Documentation, prepare presentation for client"

# Estimated / Hours Spent / New Estimate figures for rows 3-6
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = 5

$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 5

$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 5

$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 5
$ws.Range("E6").Value = 5

# Literal total for row 14 (was a SUM(D3:D13) formula)
$ws.Range("D14").Value = 20

# Move the active selection to D6
$ws.Range("D6").Select()
